$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 gets rewritten: player name/position are inserted in front of the
# existing pro-football-reference game-log columns (which shift from
# A:H to D:J), and a new numeric fantasy-points column is appended at
# the very end (L1). Only set NumberFormat="@" on the columns that would
# otherwise be auto-coerced to a date/number by plain .Value assignment
# (the date string, and the purely-numeric-looking "16" / "23.011"),
# so the rest keep their default (un-styled) General format.
$ws.Range("A1").Value = "Dunbar"
$ws.Range("B1").Value = "Steven"
$ws.Range("C1").Value = "WR"

$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "2018-12-30"

$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "16"

$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").Value = "23.011"

$ws.Range("G1").Value = "SFO"
$ws.Range("H1").Value = "@"
$ws.Range("I1").Value = "LAR"
$ws.Range("J1").Value = "L 32-48"

# K1 stays a blank placeholder cell (matches the source row's empty
# trailing column) -- just touching its format keeps the cell entry in
# the sheet without forcing a stray value into it.
$ws.Range("K1").NumberFormat = "@"

$ws.Range("L1").Value = 0
